$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: TOTAL (5-17 ans) ---
$ws.Range("C2").Value = 1212250
$ws.Range("D2").Value = 40.5
$ws.Range("E2").Value = 491129
$ws.Range("F2").Value = 4.1
$ws.Range("G2").Value = 49876
$ws.Range("H2").Value = 49.8
$ws.Range("I2").Value = 603526
$ws.Range("J2").Value = 5.6
$ws.Range("K2").Value = 67720
$ws.Range("N2").Value = 1212250

# --- Row 3: non_pdi -> hote ---
$ws.Range("A3").Value = "hote (5-17 y.o.)"
$ws.Range("B3").Value = "hote"
$ws.Range("C3").Value = 1079110
$ws.Range("D3").Value = 40.3
$ws.Range("E3").Value = 434615
$ws.Range("F3").Value = 4.1
$ws.Range("G3").Value = 43979
$ws.Range("H3").Value = 50.1
$ws.Range("I3").Value = 540693
$ws.Range("J3").Value = 5.5
$ws.Range("K3").Value = 59824
$ws.Range("N3").Value = 1079110

# --- Row 4: pdi -> idp_host ---
$ws.Range("A4").Value = "idp_host (5-17 y.o.)"
$ws.Range("B4").Value = "idp_host"
$ws.Range("C4").Value = 66749
$ws.Range("D4").Value = 44.6
$ws.Range("E4").Value = 29787
$ws.Range("F4").Value = 3.8
$ws.Range("G4").Value = 2561
$ws.Range("H4").Value = 46.3
$ws.Range("I4").Value = 30925
$ws.Range("J4").Value = 5.2
$ws.Range("K4").Value = 3476
$ws.Range("N4").Value = 66749

# --- Row 5 (new): retourne ---
$ws.Range("A5").Value = "retourne (5-17 y.o.)"
$ws.Range("B5").Value = "retourne"
$ws.Range("C5").Value = 49422
$ws.Range("D5").Value = 38.2
$ws.Range("E5").Value = 18874
$ws.Range("F5").Value = 5.4
$ws.Range("G5").Value = 2658
$ws.Range("H5").Value = 49.3
$ws.Range("I5").Value = 24380
$ws.Range("J5").Value = 7.1
$ws.Range("K5").Value = 3509
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 49422

# --- Row 6 (new): idp_site ---
$ws.Range("A6").Value = "idp_site (5-17 y.o.)"
$ws.Range("B6").Value = "idp_site"
$ws.Range("C6").Value = 16970
$ws.Range("D6").Value = 46.3
$ws.Range("E6").Value = 7853
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 678
$ws.Range("H6").Value = 44.4
$ws.Range("I6").Value = 7528
$ws.Range("J6").Value = 5.4
$ws.Range("K6").Value = 911
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 16970
